# update stanley algorithm description
#
# 1) Bump the footer date field (2023/8/9 -> 2023/8/16) everywhere it
#    appears: on the slide master and on every slide layout's
#    "Date Placeholder" shape.
# 2) On slide 4's content placeholder, split the "P4 = sf" run into
#    "P4 = " / "sf" / "（最终的弧长）" so the extra clarifying text is
#    appended after "sf".

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
                $tr = $sh.TextFrame.TextRange
                if ($tr.Text -like "*2023/8/9*") {
                    $tr.Text = "2023/8/16"
                }
            }
        }
    }
}

# --- Slide master ---
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# --- Every slide layout under the master ---
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# --- Slide 4: "P4 = sf" -> "P4 = " + "sf" + "（最终的弧长）" ---
$slide = $p.Slides.Item(4)
$shape = $slide.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$paraCount = $tr.Paragraphs().Count
$targetIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $ptext = $tr.Paragraphs($i, 1).Text
    if ($ptext.StartsWith("P4 = sf")) {
        $targetIndex = $i
    }
}

if ($targetIndex -gt 0) {
    $para = $tr.Paragraphs($targetIndex, 1)

    # Build the Chinese suffix "（最终的弧长）" via char codes (string-concat,
    # not numeric add, so start the chain from an empty string literal).
    $suffix = "" + [char]0xFF08 + [char]0x6700 + [char]0x7EC8 + [char]0x7684 + [char]0x5F27 + [char]0x957F + [char]0xFF09

    # Append the suffix right after the existing "P4 = sf" text.
    $para.Text = "P4 = sf" + $suffix

    # Re-split "P4 = sf" into two distinct runs "P4 = " and "sf" by
    # briefly mutating then restoring the "sf" substring -- identical
    # text assigned to an unchanged substring collapses back into one
    # run, but a genuine round-trip edit keeps the run boundary.
    $sfRange = $para.Characters(6, 2)
    $sfRange.Text = "zz"
    $sfRange2 = $para.Characters(6, 2)
    $sfRange2.Text = "sf"
}
